$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (I, J), using the same bold/border/
# centered formatting as the other header cells in row 1 (column H).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Per-row data for the new I (I0) and J (IF) columns.
$iVals = @{
    2 = 1;  3 = 1;  4 = 1;  5 = 1;  6 = 1;  7 = 1;  8 = 1;  9 = 1;  10 = 1;
    11 = 3; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1; 27 = 1; 28 = 1;
    29 = 1; 30 = 1; 31 = 1; 32 = 1; 33 = 1; 34 = 1; 35 = 1
}
$jVals = @{
    2 = 5;  3 = 6;  4 = 7;  5 = 5;  6 = 5;  7 = 6;  8 = 6;  9 = 8;  10 = 6;
    11 = 7; 12 = 6; 13 = 7; 14 = 6; 15 = 6; 16 = 7; 17 = 6; 18 = 7; 19 = 6;
    20 = 6; 21 = 6; 22 = 6; 23 = 6; 24 = 7; 25 = 6; 26 = 8; 27 = 6; 28 = 6;
    29 = 8; 30 = 7; 31 = 6; 32 = 1; 33 = 5; 34 = 4; 35 = 1
}

for ($row = 2; $row -le 35; $row++) {
    $ws.Cells.Item($row, 9).Value = $iVals[$row]
    $ws.Cells.Item($row, 10).Value = $jVals[$row]
}
